# Update countries & provincias Spain
# - Update "Datos actualizados" timestamp on A1
# - Austria (row 19) daily numbers updated
# - Rumania overtakes Peru/Dinamarca/Chequia (rows 32-35 shuffle down)
# - Senegal overtakes Bolivia/Estado de Palestina (rows 107-109 shuffle down)
# - Guinea-Bisau overtakes Liberia (rows 154-155 shuffle down)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 12:22"

# --- Austria (row 19) : rank unchanged, new daily figures ---
$ws.Range("B19").Value = 13730
$ws.Range("C19").Value = 170
$ws.Range("E19").Value = 6789

# --- Rumania climbs from rank 35 to rank 32 ---
# Row 32 becomes Rumania with its updated figures
$ws.Range("A32").Value = "Rumania"
$ws.Range("B32").Value = 5990
$ws.Range("C32").Value = 523
$ws.Range("D32").Value = 758
$ws.Range("E32").Value = 4950
$ws.Range("F32").Value = 208
$ws.Range("G32").Value = 12
$ws.Range("H32").Value = 282

# Row 33 becomes Peru (previously at row 32), figures unchanged
$ws.Range("A33").Value = "Peru"
$ws.Range("B33").Value = 5897
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 1569
$ws.Range("E33").Value = 4159
$ws.Range("F33").Value = 130
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 169

# Row 34 becomes Dinamarca (previously at row 33), figures unchanged
$ws.Range("A34").Value = "Dinamarca"
$ws.Range("B34").Value = 5819
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 1773
$ws.Range("E34").Value = 3799
$ws.Range("F34").Value = 113
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 247

# Row 35 becomes Chequia (previously at row 34), figures unchanged
$ws.Range("A35").Value = "Chequia"
$ws.Range("B35").Value = 5735
$ws.Range("C35").Value = 3
$ws.Range("D35").Value = 370
$ws.Range("E35").Value = 5242
$ws.Range("F35").Value = 92
$ws.Range("G35").Value = 4
$ws.Range("H35").Value = 123

# --- Senegal climbs from rank 109 to rank 107 ---
# Row 107 becomes Senegal with its updated figures
$ws.Range("A107").Value = "Senegal"
$ws.Range("B107").Value = 278
$ws.Range("C107").Value = 13
$ws.Range("D107").Value = 137
$ws.Range("E107").Value = 139
$ws.Range("F107").Value = 1
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 2

# Row 108 becomes Bolivia (previously at row 107), figures unchanged
$ws.Range("A108").Value = "Bolivia"
$ws.Range("B108").Value = 268
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 2
$ws.Range("E108").Value = 247
$ws.Range("F108").Value = 3
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 19

# Row 109 becomes Estado de Palestina (previously at row 108), figures unchanged
$ws.Range("A109").Value = "Estado de Palestina"
$ws.Range("B109").Value = 268
$ws.Range("C109").Value = 1
$ws.Range("D109").Value = 46
$ws.Range("E109").Value = 220
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 2

# --- Guinea-Bisau climbs from rank 155 to rank 154 ---
# Row 154 becomes Guinea-Bisau with its updated figures
$ws.Range("A154").Value = "Guinea-Bisau"
$ws.Range("B154").Value = 38
$ws.Range("C154").Value = 2
$ws.Range("D154").Value = 0
$ws.Range("E154").Value = 38
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 0

# Row 155 becomes Liberia (previously at row 154), figures unchanged
$ws.Range("A155").Value = "Liberia"
$ws.Range("B155").Value = 37
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 3
$ws.Range("E155").Value = 29
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 5
